$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.052.33'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.788.32'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -3.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3615'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.75'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07494'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.61'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.142'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.311'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.791.59'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.17'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001065'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06355'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.21%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.22'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.970'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.066.55'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.129'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -8.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.97'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.34'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.996.80'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.171'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.09'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.158'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.753'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08980'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.521'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.61'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02322'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.081'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6456'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2108'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06050'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.846'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.57'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5991'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.705'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.43'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.979'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.154'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06952'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.58%  '
